# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" before the existing "2022-Q2" sheet,
# populates it with the Q3 fund-holding data, and records the Q3 summary
# row on the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q3" sheet right before "2022-Q2" ----------
# Worksheets.Add() with no args inserts before the active sheet, and
# "2022-Q2" is the active sheet in this workbook.
$q3 = $wb.Worksheets.Add()
$q3.Name = "2022-Q3"

# --- 2. Fill in the "2022-Q3" sheet contents ----------------------------
# Columns B and D-G hold text (fund codes like "011346" have a leading
# zero, and the numeric-looking figures are stored as text too). Pre-format
# those cells as Text so Excel stores the values as strings instead of
# coercing them to numbers (and drops the leading zeros).
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "011346"
$q3.Range("C2").Value = "淳厚鑫淳一年持有期混合"
$q3.Range("D2").Value = "4.81"
$q3.Range("E2").Value = "69.72"
$q3.Range("F2").Value = "3.39"
$q3.Range("G2").Value = "0.1631"
$q3.Range("H2").Value = 3

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "012454"
$q3.Range("C3").Value = "淳厚鑫悦混合A"
$q3.Range("D3").Value = "2.06"
$q3.Range("E3").Value = "75.61"
$q3.Range("F3").Value = "3.43"
$q3.Range("G3").Value = "0.0707"
$q3.Range("H3").Value = 4

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "012455"
$q3.Range("C4").Value = "淳厚鑫悦混合C"
$q3.Range("D4").Value = "0.68"
$q3.Range("E4").Value = "75.61"
$q3.Range("F4").Value = "3.43"
$q3.Range("G4").Value = "0.0233"
$q3.Range("H4").Value = 4

# --- 3. Record the Q3 summary row on the totals sheet -------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 0.26

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0
